$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")

$ws.Range("Q1").Value = "person__childOf"
$ws.Range("R1").Value = "person__parentOf"
$ws.Range("S1").Value = "person__spouse"
